# Resevoir_Solution.xlsx - "Pushing stuff after exam"
# Update the adjustable-cell (decision variable) row and the unit-cost row,
# which ripple into the SUMPRODUCT totals, then move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 - adjustable cells (solver_adj range C6:H6)
$ws.Range("D6").Value = 5
$ws.Range("G6").Value = 35

# Row 7 - unit cost coefficients
$ws.Range("C7").Value = 5.4
$ws.Range("F7").Value = 2.7

# K16 is a literal RHS value (not a formula) that also needs to track the
# new supply total for Reservoir 2
$ws.Range("K16").Value = 75

# Move the active selection from G17 to F7
$ws.Range("F7").Select()
